$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.973.88"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.65%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.493.96"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.93%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "535.54"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.45%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.31"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.43%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  +0.49%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.517.89"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.15%  "

$ws.Range("E10").Value = "  -0.49%  "

$ws.Range("E11").Value = "  -2.80%  "

$ws.Range("E12").Value = "  -1.57%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.346"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -3.54%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.947.06"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.54%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.91"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.55%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "58.844.49"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.70%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000138"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.77%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.516.81"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.08%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.08"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.59%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "323.24"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.65%  "

$ws.Range("E22").Value = "  -0.04%  "

$ws.Range("E23").Value = "  +1.54%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.28"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.36%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.420"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.80%  "

$ws.Range("E26").Value = "  -1.19%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.05%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.57"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.24%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.70"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -4.10%  "

$ws.Range("E30").Value = "  -1.73%  "

$ws.Range("E31").Value = "  -1.52%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "166.46"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.29%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.16"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +4.08%  "

$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.04%  "

$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.47"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.96%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.40"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.69%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.08"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.39%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.53"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.42%  "

$ws.Range("E39").Value = "  -0.89%  "

$ws.Range("E40").Value = "  -0.43%  "

$ws.Range("E41").Value = "  -2.17%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "285.44"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.78%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.23"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.42%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "132.35"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +7.68%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.998"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.00%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.603"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.73%  "

$ws.Range("E47").Value = "  +0.13%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0925"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.32%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0507"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.60%  "

$ws.Range("E50").Value = "  -2.22%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.19"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.81%  "
